$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the athlete name in B8 from "Lê Quang Liêm" to "Co Ca"
$ws.Range("B8").Value = "Co Ca"

# Update the active selection to B8 (was C8)
$ws.Range("B8").Select()
